# feat: add 2022-Q4 data
#
# - insert a new "2022-Q4" sheet (a copy of the current "2021-Q1" sheet,
#   right after it) with the refreshed fund numbers
# - keep a "2021-Q1" sheet with the original (pre-edit) numbers
# - update the "总计" (summary) sheet with a new first data row for 2022-Q4
#   and shift the old rows down

$wb = $excel.ActiveWorkbook

$summary = $wb.Worksheets.Item(1)
$q1sheet = $wb.Worksheets.Item("2021-Q1")

# 1) Duplicate the existing "2021-Q1" worksheet; the copy lands right after
#    it and keeps all of its formatting/styles (headers, borders, etc).
$q1sheet.Copy($null, $q1sheet)
$newSheet = $wb.Worksheets.Item($q1sheet.Index + 1)

# 2) The copy becomes the new "2022-Q4" sheet; update its fund figures.
#    D2/E2/F2/G2 are text cells in the source data (not numbers), so force
#    text storage with a leading apostrophe and strip the style it adds
#    back off afterwards (keeps the cells unstyled, like the original).
$newSheet.Name = "2022-Q4"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("D2").Value = "'1.18"
$newSheet.Range("E2").Value = "'92.77"
$newSheet.Range("F2").Value = "'2.51"
$newSheet.Range("G2").Value = "'0.0296"
$newSheet.Range("D2:G2").Style = "Normal"
$newSheet.Range("H2").Value = 9

# 3) Move the new sheet so it sits between "总计" and "2021-Q1".
$newSheet.Move($null, $summary)

# 4) Update the "总计" sheet: insert a new row 2 for 2022-Q4 and push the
#    existing 2021-Q1 / 2020-Q4 rows down, renumbering column A.
$summary.Rows("2").Insert()

# match formatting of the other data rows (row 2 inherits row 1's header
# style on Insert, so re-stamp it from row 3, which still has the
# original data-row formatting)
$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)
$summary.Range("B2:D2").Style = "Normal"

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 1
$summary.Range("D2").Value = 0.03

$summary.Range("A3").Value = 1
$summary.Range("A4").Value = 2

# the "2021-Q1" copy inherited the active/selected tab from the source
# sheet it was copied from; restore the original active tab ("2020-Q4").
$q4sheet = $wb.Worksheets.Item("2020-Q4")
$q4sheet.Activate()

Write-Output "done"
